$d = $word.ActiveDocument

# Collapse the Title paragraph runs into a single merged run.
$d.Content.Find.Execute("Questions: Introduction to radians", $false, $false, $false, $false, $false, $true, 1, $false, "Questions: Introduction to radians", 2) | Out-Null

# Collapse the Author paragraph runs into a single merged run.
$d.Content.Find.Execute("Mark Toner, Ifan Howell-Baines", $false, $false, $false, $false, $false, $true, 1, $false, "Mark Toner, Ifan Howell-Baines", 2) | Out-Null

# Collapse the Abstract paragraph runs into a single merged run.
$d.Content.Find.Execute("Questions relating to the introduction to radians study guide.", $false, $false, $false, $false, $false, $true, 1, $false, "Questions relating to the introduction to radians study guide.", 2) | Out-Null
